# Apply the Dec 30 2023 cryptos-list data refresh (price/volume updates,
# plus a few rank swaps where two coins traded places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.240.33'
$ws.Range("E2").Value = '  +0.55%  '

# Row 3
$ws.Range("D3").Value = '2.291.94'
$ws.Range("E3").Value = '  +0.09%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '''316.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.65%  '

# Row 6
$ws.Range("D6").Value = '''102.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.17%  '

# Row 7
$ws.Range("E7").Value = '  +0.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.09%  '

# Row 9
$ws.Range("E9").Value = '  -0.54%  '

# Row 10
$ws.Range("D10").Value = '''39.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.10%  '

# Row 11
$ws.Range("D11").Value = '''0.0906'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.51%  '

# Row 12
$ws.Range("D12").Value = '''8.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.67%  '

# Row 13
$ws.Range("E13").Value = '  +0.62%  '

# Row 14
$ws.Range("D14").Value = '''0.957'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.70%  '

# Row 15
$ws.Range("D15").Value = '''15.21'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.68%  '

# Row 16
$ws.Range("D16").Value = '2.639.20'
$ws.Range("E16").Value = '  +0.20%  '

# Row 17
$ws.Range("D17").Value = '2.289.13'
$ws.Range("E17").Value = '  -0.03%  '

# Row 18
$ws.Range("D18").Value = '42.359.74'
$ws.Range("E18").Value = '  +1.03%  '

# Row 19
$ws.Range("E19").Value = '  -1.73%  '

# Row 20
$ws.Range("E20").Value = '  +0.73%  '

# Row 21
$ws.Range("D21").Value = '''12.26'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +32.28%  '

# Row 22
$ws.Range("D22").Value = '''73.32'
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = '''3.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.61%  '

# Row 24
$ws.Range("D24").Value = '''275.86'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +7.72%  '

# Row 25
$ws.Range("E25").Value = '  -1.95%  '

# Row 26
$ws.Range("E26").Value = '  -0.35%  '

# Row 27
$ws.Range("D27").Value = '''10.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.33%  '

# Row 28
$ws.Range("D28").Value = '''2.36'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.66%  '

# Row 29
$ws.Range("D29").Value = '''22.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.12%  '

# Row 30
$ws.Range("D30").Value = '''37.48'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.66%  '

# Row 31
$ws.Range("D31").Value = '''165.96'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.51%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''6.02'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.77%  '

# Row 33
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '''0.0872'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.01%  '

# Row 34
$ws.Range("E34").Value = '  +3.37%  '

# Row 35
$ws.Range("D35").Value = '''2.67'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.22%  '

# Row 36
$ws.Range("E36").Value = '  -0.85%  '

# Row 37
$ws.Range("E37").Value = '  -0.23%  '

# Row 38
$ws.Range("D38").Value = '''0.0361'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.33%  '

# Row 39
$ws.Range("D39").Value = '''3.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.74%  '

# Row 40
$ws.Range("D40").Value = '''2.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.33%  '

# Row 41
$ws.Range("D41").Value = '''1.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.13%  '

# Row 42
$ws.Range("B42").Value = 'BitcoinSV'
$ws.Range("C42").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D42").Value = '''96.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.74%  '

# Row 43
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").Value = '''69.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.18%  '

# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.09%  '

# Row 45
$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = '''0.225'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.10%  '

# Row 46
$ws.Range("D46").Value = '''11.92'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.27%  '

# Row 47
$ws.Range("D47").Value = '''112.88'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.70%  '

# Row 48
$ws.Range("D48").Value = '''79.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.22%  '

# Row 49
$ws.Range("D49").Value = '''8.95'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.86%  '

# Row 50
$ws.Range("E50").Value = '  -0.52%  '

# Row 51
$ws.Range("D51").Value = '1.596.33'
$ws.Range("E51").Value = '  +3.01%  '
